$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "campus"
$ws.Range("H2").Value = "B"

$ws.Range("H2").Select()
